# Add a new worksheet "No Header Works" as the last sheet in the workbook,
# populate it with two rows of header-less data, and leave it as the active sheet.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "No Header Works"

$newSheet.Range("A1").Value = "Value1"
$newSheet.Range("B1").Value = 1234
$newSheet.Range("A2").Value = "Value2"
$newSheet.Range("B2").Value = 3456

$newSheet.Range("B2").Select() | Out-Null
